# ABr // Documentation/Report // Work on the Report
#
# Fix a typo in the "Controller Structure" bullet point description:
# "Forumlas" -> "Formulas"

$d = $word.ActiveDocument

$d.Content.Find.Execute("Forumlas", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Formulas", 2)
